$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the XML diff (values stored as text to match
# the workbook convention of inlineStr / shared-string text cells).
$textFormatCells = @(
    "D2",
    "D3",
    "D4",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D20",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D47",
    "D48",
    "D49",
    "D50"
)

foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "245.82"
$ws.Range("D3").Value = "25.37"
$ws.Range("D4").Value = "5.043"
$ws.Range("D6").Value = "6.567"
$ws.Range("D7").Value = "3.011"
$ws.Range("D8").Value = "0.8184"
$ws.Range("D9").Value = "0.8350"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.009692"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1338"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.06955"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("D13").Value = "0.02825"
$ws.Range("D14").Value = "0.09403"
$ws.Range("D15").Value = "0.001521"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006177"
$ws.Range("E16").Value = "15TigerCashTCH"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.496"
$ws.Range("E17").Value = "16LEOLEO"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "2.092"
$ws.Range("E18").Value = "17BTSETokenBTSE"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3188"
$ws.Range("E19").Value = "18BitpandaEcosystemTokenBEST"
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "0.03224"
$ws.Range("E20").Value = "19LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("D22").Value = "3.740"
$ws.Range("D23").Value = "0.04695"
$ws.Range("D24").Value = "0.1341"
$ws.Range("D25").Value = "0.001243"
$ws.Range("D26").Value = "0.004291"
$ws.Range("D27").Value = "0.00009699"
$ws.Range("E27").Value = "26NitroExNTX"
$ws.Range("D28").Value = "0.0001940"
$ws.Range("D40").Value = "0.03661"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006219"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1055"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("D43").Value = "0.002600"
$ws.Range("D44").Value = "0.007432"
$ws.Range("D45").Value = "0.00005293"
$ws.Range("D47").Value = "0.1800"
$ws.Range("D48").Value = "0.002016"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("D50").Value = "0.0002000"

foreach ($addr in $textFormatCells) {
    $ws.Range($addr).ClearFormats()
}
